$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3816.2
$ws.Range("I76").Value = 3318.2222
$ws.Range("K76").Value = 3318.2222
$ws.Range("M76").Value = -3003.2222
$ws.Range("H79").Value = 3816.2
$ws.Range("I79").Value = 3318.2222
$ws.Range("K79").Value = 3318.2222
$ws.Range("M79").Value = -2226.2222
$ws.Range("H132").Value = 3269.4167
$ws.Range("I132").Value = 3198.4888
$ws.Range("K132").Value = 9595.466400000001
$ws.Range("M132").Value = -7065.466400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 109666
$ws.Range("I28").Value = 100001
$ws.Range("J28").Value = 114498.5
$ws.Range("K28").Value = 100001
$ws.Range("L28").Value = 114498.5
$ws.Range("M28").Value = -99809
$ws.Range("N28").Value = -114882.5
$ws.Range("H32").Value = 10419579
$ws.Range("I32").Value = 6411952.5
$ws.Range("K32").Value = 6411952.5
$ws.Range("M32").Value = -6411665.5
$ws.Range("H61").Value = 1967.6731
$ws.Range("I61").Value = 1826.4
$ws.Range("J61").Value = 5499.5
$ws.Range("K61").Value = 1826.4
$ws.Range("L61").Value = 5499.5
$ws.Range("M61").Value = -1614.4
$ws.Range("N61").Value = -5923.5
$ws.Range("H63").Value = 3526.9092
$ws.Range("I63").Value = 2428
$ws.Range("K63").Value = 2428
$ws.Range("M63").Value = -1742
$ws.Range("H66").Value = 3526.9092
$ws.Range("I66").Value = 2428
$ws.Range("K66").Value = 12140
$ws.Range("M66").Value = -8708
$ws.Range("H74").Value = 1220.3704
$ws.Range("I74").Value = 1009.6539
$ws.Range("K74").Value = 1009.6539
$ws.Range("M74").Value = -135.6539
$ws.Range("H77").Value = 1220.3704
$ws.Range("I77").Value = 1009.6539
$ws.Range("K77").Value = 5048.2695
$ws.Range("M77").Value = -680.2695000000003
$ws.Range("H99").Value = 109666
$ws.Range("I99").Value = 100001
$ws.Range("J99").Value = 114498.5
$ws.Range("K99").Value = 100001
$ws.Range("L99").Value = 114498.5
$ws.Range("M99").Value = -97006
$ws.Range("N99").Value = -120488.5
$ws.Range("H102").Value = 1427.5
$ws.Range("I102").Value = 625.6667
$ws.Range("J102").Value = 3833
$ws.Range("K102").Value = 625.6667
$ws.Range("L102").Value = 3833
$ws.Range("M102").Value = 996.3333
$ws.Range("N102").Value = -7077
$ws.Range("H122").Value = 3883.353
$ws.Range("I122").Value = 2866.7114
$ws.Range("K122").Value = 8600.1342
$ws.Range("M122").Value = -6150.1342
$ws.Range("H128").Value = 149914.5
$ws.Range("J128").Value = 149914.5
$ws.Range("L128").Value = 149914.5
$ws.Range("N128").Value = -159874.5
$ws.Range("H132").Value = 3677.6191
$ws.Range("I132").Value = 3313.3333
$ws.Range("J132").Value = 4163.3335
$ws.Range("K132").Value = 9939.999899999999
$ws.Range("L132").Value = 12490.0005
$ws.Range("M132").Value = -7409.999899999999
$ws.Range("N132").Value = -17550.0005
$ws.Range("H136").Value = 1967.6731
$ws.Range("I136").Value = 1826.4
$ws.Range("J136").Value = 5499.5
$ws.Range("K136").Value = 5479.200000000001
$ws.Range("L136").Value = 16498.5
$ws.Range("M136").Value = -2929.200000000001
$ws.Range("N136").Value = -21598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3342.0667
$ws.Range("I105").Value = 1676.4
$ws.Range("K105").Value = 1676.4
$ws.Range("M105").Value = 70.59999999999991
$ws.Range("H134").Value = 10941521
$ws.Range("I134").Value = 2233779.5
$ws.Range("J134").Value = 66671068
$ws.Range("K134").Value = 6701338.5
$ws.Range("L134").Value = 200013204
$ws.Range("M134").Value = -6698803.5
$ws.Range("N134").Value = -200018274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3944.4075
$ws.Range("J31").Value = 6218.154
$ws.Range("L31").Value = 6218.154
$ws.Range("N31").Value = -6808.154
$ws.Range("H34").Value = 3944.4075
$ws.Range("J34").Value = 6218.154
$ws.Range("L34").Value = 6218.154
$ws.Range("N34").Value = -6622.154
$ws.Range("H58").Value = 2485.5278
$ws.Range("I58").Value = 1919.1786
$ws.Range("K58").Value = 1919.1786
$ws.Range("M58").Value = -1716.1786
$ws.Range("H99").Value = 2242.5
$ws.Range("I99").Value = 2185
$ws.Range("K99").Value = 2185
$ws.Range("M99").Value = -687
$ws.Range("H126").Value = 2242.5
$ws.Range("I126").Value = 2185
$ws.Range("K126").Value = 6555
$ws.Range("M126").Value = -4085
$ws.Range("H132").Value = 2205.05
$ws.Range("I132").Value = 2182.4119
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 6547.2357
$ws.Range("L132").Value = 6999.999899999999
$ws.Range("M132").Value = -4017.2357
$ws.Range("N132").Value = -12059.9999
$ws.Range("H136").Value = 2485.5278
$ws.Range("I136").Value = 1919.1786
$ws.Range("K136").Value = 5757.5358
$ws.Range("M136").Value = -3207.5358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 17.571428
$ws.Range("I8").Value = 17.571428
$ws.Range("K8").Value = 52.71428400000001
$ws.Range("M8").Value = 86.28571599999999
$ws.Range("H132").Value = 1557.6666
$ws.Range("I132").Value = 979.375
$ws.Range("K132").Value = 8814.375
$ws.Range("M132").Value = -6284.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2632
$ws.Range("I132").Value = 2459.5881
$ws.Range("J132").Value = 4097.5
$ws.Range("K132").Value = 7378.7643
$ws.Range("L132").Value = 12292.5
$ws.Range("M132").Value = -4848.7643
$ws.Range("N132").Value = -17352.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4348.154
$ws.Range("I7").Value = 3250.1428
$ws.Range("K7").Value = 3250.1428
$ws.Range("M7").Value = -3138.1428
$ws.Range("H40").Value = 5562.7036
$ws.Range("I40").Value = 5696.4165
$ws.Range("K40").Value = 5696.4165
$ws.Range("M40").Value = -5560.4165
$ws.Range("H61").Value = 1828.1428
$ws.Range("I61").Value = 1799.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1799.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1597.5
$ws.Range("N61").Value = -2404
$ws.Range("H113").Value = 1828.1428
$ws.Range("I113").Value = 1799.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1799.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 370.5
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 3677.5557
$ws.Range("I122").Value = 3848.3076
$ws.Range("K122").Value = 11544.9228
$ws.Range("M122").Value = -9094.9228
$ws.Range("H126").Value = 4348.154
$ws.Range("I126").Value = 3250.1428
$ws.Range("K126").Value = 9750.428400000001
$ws.Range("M126").Value = -7280.428400000001
$ws.Range("H131").Value = 75000
$ws.Range("I131").Value = 50000
$ws.Range("J131").Value = 100000
$ws.Range("K131").Value = 50000
$ws.Range("L131").Value = 100000
$ws.Range("M131").Value = -44960
$ws.Range("N131").Value = -110080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 625795
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -887
$ws.Range("H27").Value = 33525.5
$ws.Range("J27").Value = 33525.5
$ws.Range("L27").Value = 33525.5
$ws.Range("N27").Value = -33663.5
$ws.Range("H62").Value = 4835.4707
$ws.Range("I62").Value = 4267.5
$ws.Range("K62").Value = 4267.5
$ws.Range("M62").Value = -3643.5
$ws.Range("H65").Value = 4835.4707
$ws.Range("I65").Value = 4267.5
$ws.Range("K65").Value = 21337.5
$ws.Range("M65").Value = -18217.5
$ws.Range("H107").Value = 605.25
$ws.Range("I107").Value = 525.5
$ws.Range("J107").Value = 685
$ws.Range("K107").Value = 1576.5
$ws.Range("L107").Value = 2055
$ws.Range("M107").Value = 343.5
$ws.Range("N107").Value = -5895
$ws.Range("H115").Value = 104964.5
$ws.Range("J115").Value = 104964.5
$ws.Range("L115").Value = 104964.5
$ws.Range("N115").Value = -108098.5
$ws.Range("H122").Value = 2177.16
$ws.Range("I122").Value = 1720.091
$ws.Range("K122").Value = 5160.272999999999
$ws.Range("M122").Value = -2710.272999999999
